$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 becomes the new "Poder Judicial de Honduras" entry
$ws.Range("A2").Value = "Poder Judicial de Honduras"

$ws.Range("D2").Value = "La potestad de impartir justicia emana del pueblo y se imparte gratuitamente en nombre del Estado, por magistrados y jueces independientes, únicamente sometidos a la Constitución y a las leyes. El Poder Judicial se integra por una Corte Suprema de Justicia, por las Cortes de Apelaciones, los Juzgados, y demás dependencias que señale la Ley"

$ws.Range("F2").Value = "Se autoriza para que todo Servidor Judicial mayor a 60 años de edad, permanezca en su`ncasa y que, dependiendo de la naturaleza de las tareas asignadas, el trabajo bajo su`nresponsabilidad pueda desarrollarlo de manera remota a través de las plataformas`ninformáticas disponibles o de la forma como sea indicada por su superior jerárquico."

# G2 previously linked to the Chile "Dirección del trabajo" site; drop that
# hyperlink before re-pointing the cell at the Honduras source.
$ws.Range("G2").Hyperlinks.Delete()
$ws.Range("G2").Value = "https://covid19honduras.org/?q=comunicado-del-poder-judicial"
$ws.Hyperlinks.Add($ws.Range("G2"), "https://covid19honduras.org/?q=comunicado-del-poder-judicial")

$ws.Range("E2").Value = "https://covid19honduras.org/?q=comunicado-del-poder-judicial"
$ws.Hyperlinks.Add($ws.Range("E2"), "https://covid19honduras.org/?q=comunicado-del-poder-judicial")

$ws.Range("H2").ClearContents()

$ws.Range("I2").Value = "15/3/2020"
$ws.Range("J2").Value = "Honduras"
$ws.Range("K2").ClearContents()

# Clear the now-unused "Descripción información" cells for the Chile rows
$ws.Range("F3").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("F5").ClearContents()
